$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1) onto
# the three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-60) gets the season win/loss/tie record appended as
# plain numeric values in the new AD/AE/AF columns.
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 30).Value = 54
    $ws.Cells.Item($r, 31).Value = 108
    $ws.Cells.Item($r, 32).Value = 0
}
